$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)

# Row 5: MIKU PARTY entry removed; old row6 content (小马宝莉only) moves up with updates
$ws1.Range("B5").Value = "2024-07-28"
$ws1.Range("C5").Value = "广州·小马宝莉only（取消）"
$ws1.Range("D5").Value = "鸿盛二路巨大创意产业园 巨大产业园·智汇港"
$ws1.Range("E5").Value = "2024.07.28 10:00-07.28 17:00"
$ws1.Range("F5").Value = 340
$ws1.Range("G5").Value = "不可售"
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=88110"
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202406/XH2hCwrg1719393458914.jpeg"

# Row 6: old row7 content (运动番only7.0) moves up with updates
$ws1.Range("C6").Value = "广州·运动番only7.0"
$ws1.Range("D6").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Range("F6").Value = 1241
$ws1.Range("G6").Value = 60
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=88473"
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202407/e98msNdA1721383295549.jpeg"

# Row 7: new entry (Le plaisir...)
$ws1.Range("B7").Value = "2024-08-03"
$ws1.Range("C7").Value = "广州·Le plaisir 第五人格&明日方舟主题同人派对（取消）"
$ws1.Range("D7").Value = "太和岗路18号负一层 8+1 live house"
$ws1.Range("E7").Value = "2024.08.03 13:00-08.03 20:00"
$ws1.Range("F7").Value = 87
$ws1.Range("G7").Value = "不可售"
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88654"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202407/B3gUl2Gn1720073290274.jpeg"

# F-column (want-to-go count) updates for unaffected rows
$ws1.Range("F3").Value = 11446
$ws1.Range("F4").Value = 1306
$ws1.Range("F8").Value = 177
$ws1.Range("F9").Value = 967
$ws1.Range("F11").Value = 2357
$ws1.Range("F13").Value = 1160
$ws1.Range("F14").Value = 881
$ws1.Range("F15").Value = 585
$ws1.Range("F17").Value = 1050
$ws1.Range("F19").Value = 108
$ws1.Range("F20").Value = 701
$ws1.Range("F21").Value = 732
$ws1.Range("F22").Value = 161
$ws1.Range("F23").Value = 430
$ws1.Range("F24").Value = 1081
$ws1.Range("F26").Value = 504
$ws1.Range("F27").Value = 553
$ws1.Range("F29").Value = 280
$ws1.Range("F30").Value = 279
$ws1.Range("F31").Value = 649
$ws1.Range("F32").Value = 2825
$ws1.Range("F33").Value = 448
$ws1.Range("F37").Value = 102
$ws1.Range("F38").Value = 1540
$ws1.Range("F41").Value = 70
$ws1.Range("F42").Value = 122
$ws1.Range("F43").Value = 57
$ws1.Range("F45").Value = 98
$ws1.Range("F46").Value = 109
$ws1.Range("F47").Value = 70

# G-column updates
$ws1.Range("G32").Value = 75

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 79
$ws2.Range("F7").Value = 96
$ws2.Range("F10").Value = 164
$ws2.Range("F11").Value = 4411
$ws2.Range("F14").Value = 137

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 2244
$ws3.Range("F3").Value = 696
$ws3.Range("F4").Value = 657

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 2244
$ws4.Range("F3").Value = 696
$ws4.Range("F5").Value = 11446
$ws4.Range("F6").Value = 657
$ws4.Range("F7").Value = 1241
$ws4.Range("F9").Value = 79
$ws4.Range("F10").Value = 967
$ws4.Range("F11").Value = 2357
$ws4.Range("F13").Value = 1160
$ws4.Range("F14").Value = 881
$ws4.Range("F15").Value = 585
$ws4.Range("F17").Value = 1050
$ws4.Range("F20").Value = 108
$ws4.Range("F21").Value = 701
$ws4.Range("F23").Value = 732
$ws4.Range("F24").Value = 161
$ws4.Range("F25").Value = 430
$ws4.Range("F26").Value = 1081
$ws4.Range("F27").Value = 96
$ws4.Range("F29").Value = 504
$ws4.Range("F30").Value = 553
$ws4.Range("F32").Value = 279
$ws4.Range("F33").Value = 2825
$ws4.Range("F34").Value = 164
$ws4.Range("F35").Value = 448
$ws4.Range("F38").Value = 102
$ws4.Range("F39").Value = 1540
$ws4.Range("F42").Value = 122
$ws4.Range("F43").Value = 57
$ws4.Range("F45").Value = 98
$ws4.Range("F47").Value = 70

$ws4.Range("G33").Value = 75
